# tugas5.xlsx: tambah kolom "Id" di Tabel Detail Pesanan (3NF) dan
# tambah diagram "Tabel Produk" kedua di bagian bawah sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Kolom "Id" baru di tabel 3NF (M36:M40), disisipkan di antara
#        "Tabel Detail Pesanan" (G:K) dan "Tabel Harga" (N:O) ---

# M36: sel kosong dengan garis bawah tipis saja (seperti header lain di atas tabel)
$ws.Range("M36").Borders.Item(9).LineStyle = 1
$ws.Range("M36").Borders.Item(9).Weight = 2

# M37: header "Id" - salin format dari header tabel yang sudah ada (style 1)
$ws.Range("K37").Copy()
$ws.Range("M37").PasteSpecial(-4122)
$ws.Range("M37").Value = "Id"

# M38:M40: data 1,2,3 - salin format dari kolom data yang sudah ada (style 2)
$ws.Range("G38").Copy()
$ws.Range("M38").PasteSpecial(-4122)
$ws.Range("M38").Value = 1

$ws.Range("G39").Copy()
$ws.Range("M39").PasteSpecial(-4122)
$ws.Range("M39").Value = 2

$ws.Range("G40").Copy()
$ws.Range("M40").PasteSpecial(-4122)
$ws.Range("M40").Value = 3

# --- 2. Diagram kedua "Tabel Produk" (G51:H55), seperti yang sudah
#        ada di G20:I25 ---

# G51: judul tabel (teks polos, tanpa style), sama dengan "Tabel Produk" di G20
$ws.Range("G51").Value = "Tabel Produk"

# G52/H52: header kolom "Produk_Id" / "Kategori" - salin format header (style 1)
$ws.Range("G21").Copy()
$ws.Range("G52").PasteSpecial(-4122)
$ws.Range("G52").Value = "Produk_Id"

$ws.Range("H21").Copy()
$ws.Range("H52").PasteSpecial(-4122)
$ws.Range("H52").Value = "Kategori"

# G53:H55: data - salin format data (style 2)
$ws.Range("G22").Copy()
$ws.Range("G53").PasteSpecial(-4122)
$ws.Range("G53").Value = 1

$ws.Range("H22").Copy()
$ws.Range("H53").PasteSpecial(-4122)
$ws.Range("H53").Value = "Elektronik"

$ws.Range("G23").Copy()
$ws.Range("G54").PasteSpecial(-4122)
$ws.Range("G54").Value = 2

$ws.Range("H23").Copy()
$ws.Range("H54").PasteSpecial(-4122)
$ws.Range("H54").Value = "Makanan"

$ws.Range("G24").Copy()
$ws.Range("G55").PasteSpecial(-4122)
$ws.Range("G55").Value = 3

$ws.Range("H24").Copy()
$ws.Range("H55").PasteSpecial(-4122)
$ws.Range("H55").Value = "Minuman"

# --- 3. Sheet view: scroll ke bawah, zoom 85%, pilih H52 ---
$win = $excel.ActiveWindow
$win.Zoom = 85
$ws.Range("H52").Select()
